# Prelim work on localization
# Adds a block of new "Electricity / E*" process & commodity set rows to the
# "SetsEditor- Proc" sheet, widens column E to fit the new (longer) entries,
# and updates the active-sheet / selection state left behind by the edit
# session (SetsEditor- Proc becomes the active tab; VEDA_Sets-Comm is no
# longer tabSelected; VEDA_Sets-Proc's lingering full-column selection is
# cleared to a normal single-cell selection).

$wb = $excel.ActiveWorkbook

$wsProc  = $wb.Worksheets.Item("SetsEditor- Proc")
$wsComm  = $wb.Worksheets.Item("VEDA_Sets-Comm")
$wsVProc = $wb.Worksheets.Item("VEDA_Sets-Proc")

# --- New data rows (12-31) on "SetsEditor- Proc" -----------------------

$rows = @(
    @{ Row = 12; E = "EPumpStorage";      F = "EPumpStorage";      H = "EPTSTO*" },
    @{ Row = 13; E = "EBiomass";          F = "EBiomass";          H = "ERB*" },
    @{ Row = 14; E = "EHydro";            F = "EHydro";            H = "ERHYD*,-ERHYD*-I" },
    @{ Row = 15; E = "EImports";          F = "EImports";          H = "ERHYD*-I" },
    @{ Row = 16; E = "EPV_Grid";          F = "EPV_Grid";          H = "ERSOLPC*" },
    @{ Row = 17; E = "EPV_RfT";           F = "EPV_RfTpIND";       H = "ERSOLPI*,ERSOLPRC*,ERSOLPRR*" },
    @{ Row = 18; E = "ECSP";              F = "ECSP";              H = "ERSOLT*" },
    @{ Row = 19; E = "EWind";             F = "EWind";             H = "ERWND*" },
    @{ Row = 20; E = "EBattery-Dist";     F = "EBattery-Commerce"; H = "ESTSC*,ESTSI*,ESTSR*" },
    @{ Row = 21; E = "EBattery-Utility";  F = "EBattery-Utility";  H = "ESTSU*" },
    @{ Row = 22; E = "ECoal";             F = "ECoal";             H = "ETC*" },
    @{ Row = 23; E = "EGas";              F = "EGas";              H = "ETG*" },
    @{ Row = 24; E = "ENuclear";          F = "ENuclear";          H = "ETN*" },
    @{ Row = 25; E = "EOil";              F = "EOil";              H = "ETO*" },
    @{ Row = 26; E = "Etrans";            F = "Etrans";            H = "ETRANS*" },
    @{ Row = 27; E = "Edist_Residential"; F = "Edist_Residential"; H = "XRESELC" },
    @{ Row = 28; E = "Edist_Agriculture"; F = "Edist_Agriculture"; H = "XAGRELC" },
    @{ Row = 29; E = "Edist_Industry";    F = "Edist_Industry";    H = "XI*ELC,XU*ELC" },
    @{ Row = 30; E = "Edist_Commercial";  F = "Edist_Commercial";  H = "XCOMELC" },
    @{ Row = 31; E = "Edist_Transport";   F = "Edist_Transport";   H = "XTRAELC" }
)

foreach ($r in $rows) {
    $wsProc.Cells.Item($r.Row, 5).Value = $r.E   # column E
    $wsProc.Cells.Item($r.Row, 6).Value = $r.F   # column F
    $wsProc.Cells.Item($r.Row, 8).Value = $r.H   # column H
}

# --- Column E is now wider (it holds the longer new entries) -----------

$wsProc.Columns.Item(5).ColumnWidth = 18.140625

# --- Active sheet / selection bookkeeping -------------------------------

# Tidy up the leftover "select everything" selection on VEDA_Sets-Proc.
[void]$wsVProc.Range("D3").Select()

# SetsEditor- Proc becomes the active sheet, with the freshly-entered block
# selected.
[void]$wsProc.Activate()
[void]$wsProc.Range("A18:XFD19").Select()
